$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk update column F (Data) for rows 2-51: 17-12-2022 -> 18-12-2022
$ws.Range("F2:F51").Value = "18-12-2022"

# Bulk update column G (Hora) for rows 2-51: 23 -> 0 (force text to match original inline-string formatting)
$ws.Range("G2:G51").Value = "'0"

# Per-cell updates for column D (Price)
$ws.Range("D2").Value = "'240.32"
$ws.Range("D3").Value = "'21.96"
$ws.Range("D4").Value = "'5.450"
$ws.Range("D5").Value = "'0.05596"
$ws.Range("D6").Value = "'6.479"
$ws.Range("D7").Value = "'3.365"
$ws.Range("D8").Value = "'0.8065"
$ws.Range("D9").Value = "'1.071"
$ws.Range("D10").Value = "'0.1426"
$ws.Range("D11").Value = "'0.07380"
$ws.Range("D12").Value = "'0.03272"
$ws.Range("D13").Value = "'0.02941"
$ws.Range("D14").Value = "'0.09240"
$ws.Range("D15").Value = "'0.001664"
$ws.Range("D16").Value = "'3.251"
$ws.Range("D17").Value = "'0.04804"
$ws.Range("D18").Value = "'0.0005744"
$ws.Range("D19").Value = "'0.006253"
$ws.Range("D20").Value = "'0.001046"
$ws.Range("D21").Value = "'0.003798"
$ws.Range("D22").Value = "'0.0001498"
$ws.Range("D23").Value = "'0.0004795"
$ws.Range("D24").Value = "'3.974"
$ws.Range("D25").Value = "'2.198"
$ws.Range("D27").Value = "'0.1312"
$ws.Range("D40").Value = "'0.04196"
$ws.Range("D41").Value = "'0.006977"
$ws.Range("D42").Value = "'0.003496"
$ws.Range("D43").Value = "'0.1043"
$ws.Range("D44").Value = "'0.008807"
$ws.Range("D45").Value = "'0.00005461"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.6792"
$ws.Range("D48").Value = "'0.03064"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D50").Value = "'0.01009"

# Per-cell updates for column E (Volume(1h))
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
